$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert three new columns at the front (A:C), shifting the existing
# question-bank columns (old A..K) to the right (new D..N).
$ws.Range("A:C").Insert()

# Populate the three new columns with headers + sample row values.
# Written in this order (A, C, B) so the shared-string table picks up
# academic_id / course_id / sem_id in that sequence.
$ws.Range("A1").Value = "academic_id"
$ws.Range("C1").Value = "course_id"
$ws.Range("B1").Value = "sem_id"

$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 17
$ws.Range("B2").Value = 1

# Widen the two new id columns, matching the author's manual resize.
$ws.Range("A:B").ColumnWidth = 12.166666666666666

# Approximate the author's manual column widths on the shifted-right
# columns (closest values reachable through this host's width stepping).
$ws.Range("G:G").ColumnWidth = 38.5
$ws.Range("H:H").ColumnWidth = 10.5
$ws.Range("I:I").ColumnWidth = 16.333333333333336
$ws.Range("J:J").ColumnWidth = 15.499999999999998
$ws.Range("L:L").ColumnWidth = 17.833333333333336
$ws.Range("M:M").ColumnWidth = 12.666666666666666
$ws.Range("N:N").ColumnWidth = 28.666666666666668

# Move the active selection, matching the author's final cursor position.
$ws.Range("E8").Select() | Out-Null

# The question_type / filter-database helper range moved down one column
# (N -> O) when the author separately inserted a column ahead of it.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$O`$2:`$O`$3"
$wb.Names.Item("question_type").RefersTo = "=Sheet1!`$O`$2:`$O`$3"
$wb.Names.Item("type").RefersTo = "=Sheet1!`$O`$2:`$O`$6"
